# 2022April Sprint Test cases.xlsx -- "Update 2022April Sprint Test"
#
# Inserts a new column G ("Lookup information" detail notes) in front of the
# existing Notes column (which shifts to H), populates it, tweaks a handful
# of row heights / wrap settings, widens the new column, and resets the
# active selection + page orientation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert a new column before the existing "Notes" column (G -> H).
#    Excel copies column F's format (wrap text, left/vcenter) onto the
#    freshly inserted column G automatically.
# ---------------------------------------------------------------------
$ws.Columns("G:G").Insert()

# ---------------------------------------------------------------------
# 2. Populate the new column G with "Lookup information" detail text.
#    Rows whose text should stay short (AppliedPatch_* / plain notes)
#    reuse exactly the same string that is already sitting in H;
#    rows 4-7 get brand-new, longer "Assigned Queue does not exist..."
#    explanations (and their H neighbour also changes for rows 4-7).
#    Write order here follows the original authoring order so new
#    shared-string entries land at the same table indices.
# ---------------------------------------------------------------------
$ws.Range("G2").Value = "not added to queue by Dispatcher"
$ws.Range("G3").Value = "Reason: No matching Account Number on the lookup table. Account Number: TESTING"

$ws.Range("H5").Value = "Reason: SR Failed to update. Lookup InformationSR0003015453,Beverly Johnson,COR-Bad-Queue-Name,991047.."
$ws.Range("H6").Value = "Reason: SR Failed to update. Lookup Information:SR0003015454,COR-Bad-Queue-Name,991048,Beverly Johnson."
$ws.Range("H7").Value = "Reason: SR Failed to update. Lookup InformationSR0003015455,Beverly Johnson,COR-Bad-Queue-Name,991040.,Beverly Johnson."

$ws.Range("G4").Value = "Reason: SR Failed to update. Assigned Queue does not exist in OEC. Lookup information: SR0003015452, COR-Bad-Queue-Name."
$ws.Range("H4").Value = "Reason: SR Failed to update. Assigned Queue does not exist in OEC. Lookup information: SR0003015452, COR-Bad-Queue-Name."

$ws.Range("G5").Value = "Reason: SR Failed to update. Assigned Queue does not exist in OEC. Lookup information: SR0003015453, Beverly Johnson, COR-Bad-Queue-Name, 991047. ."
$ws.Range("G6").Value = "Reason: SR Failed to update. Assigned Queue does not exist in OEC. Lookup information: SR0003015454, COR-Bad-Queue-Name, 991048, Beverly Johnson."
$ws.Range("G7").Value = "Reason: SR Failed to update. Assigned Queue does not exist in OEC. Lookup information: SR0003015455, Beverly Johnson, COR-Bad-Queue-Name, 991040. , Beverly Johnson."

$ws.Range("G8").Value = "AppliedPatch_OnlyQueue"
$ws.Range("G10").Value = "AppliedPatch_AssignTo"
$ws.Range("G12").Value = "AppliedPatch_CCRAssignTo"
# G9, G11, G13, G14 stay blank -- their values remain only in column H.

# ---------------------------------------------------------------------
# 3. Column G's inserted format wraps text everywhere (copied from F).
#    The short "note" rows should NOT wrap (matches the look of column H),
#    so turn wrapping back off for those.
# ---------------------------------------------------------------------
$ws.Range("G2").WrapText = $false
$ws.Range("G8").WrapText = $false
$ws.Range("G9").WrapText = $false
$ws.Range("G10").WrapText = $false
$ws.Range("G11").WrapText = $false
$ws.Range("G12").WrapText = $false

# ---------------------------------------------------------------------
# 4. Row heights grow to fit the new multi-line detail text.
# ---------------------------------------------------------------------
$ws.Rows(4).RowHeight = 60
$ws.Rows(5).RowHeight = 75
$ws.Rows(6).RowHeight = 75
$ws.Rows(7).RowHeight = 75

# ---------------------------------------------------------------------
# 5. Widen the new column G and re-select the cell the author ended on.
# ---------------------------------------------------------------------
$ws.Columns("G").ColumnWidth = 34.67

$ws.Range("E6").Select()

# ---------------------------------------------------------------------
# 6. Page orientation was (re)set to portrait.
# ---------------------------------------------------------------------
$ws.PageSetup.Orientation = 1
